$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values keep exact
# string representation (trailing zeros, precision) instead of being coerced
# to numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.153.53'
$ws.Range("D3").Value = '1.750.82'
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '236.75'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.5342'
$ws.Range("E7").Value = '  +2.29%  '
$ws.Range("D8").Value = '0.2810'
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("D9").Value = '0.06178'
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").Value = '1.745.89'
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").Value = '0.07180'
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = '0.6486'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '4.632'
$ws.Range("E14").Value = '  +2.09%  '
$ws.Range("D15").Value = '78.53'
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("D16").Value = '0.9997'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '0.9997'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '26.040.16'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '11.78'
$ws.Range("E19").Value = '  +2.27%  '
$ws.Range("D20").Value = '0.000006763'
$ws.Range("E20").Value = '  +2.00%  '
$ws.Range("D21").Value = '1.968.45'
$ws.Range("D22").Value = '4.350'
$ws.Range("E22").Value = '  +4.28%  '
$ws.Range("D23").Value = '8.728'
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("D24").Value = '5.249'
$ws.Range("E24").Value = '  +1.69%  '
$ws.Range("D25").Value = '139.36'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '1.522'
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("D28").Value = '1.801'
$ws.Range("E28").Value = '  -2.27%  '
$ws.Range("D29").Value = '105.17'
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("D30").Value = '0.08325'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '3.778'
$ws.Range("E31").Value = '  +3.05%  '
$ws.Range("E32").Value = '  +5.70%  '
$ws.Range("D33").Value = '0.04645'
$ws.Range("E33").Value = '  +3.77%  '
$ws.Range("D34").Value = '2.645'
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").Value = '1.010'
$ws.Range("E35").Value = '  +2.28%  '
$ws.Range("D36").Value = '0.6322'
$ws.Range("E36").Value = '  +3.18%  '
$ws.Range("D37").Value = '2.711'
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("D39").Value = '1.967'
$ws.Range("D40").Value = '0.9995'
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").Value = '102.20'
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("D42").Value = '0.3932'
$ws.Range("E42").Value = '  +1.38%  '
$ws.Range("D43").Value = '0.7515'
$ws.Range("E43").Value = '  +2.04%  '
$ws.Range("D44").Value = '5.077'
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").Value = '0.1155'
$ws.Range("E45").Value = '  +3.16%  '
$ws.Range("D46").Value = '6.366'
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("D47").Value = '0.05340'
$ws.Range("E47").Value = '  -2.49%  '
$ws.Range("D48").Value = '54.75'
$ws.Range("E48").Value = '  +3.19%  '
$ws.Range("D49").Value = '31.01'
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("D50").Value = '0.3488'
$ws.Range("E50").Value = '  +1.23%  '
$ws.Range("D51").Value = '7.636'
$ws.Range("E51").Value = '  +0.26%  '
